# Virtualbox Guest Additions added for Debian
# Adds a new "virtualbox-guest-additions" row to the support matrix and
# tweaks a handful of existing Supported/Not-supported cell colors.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Stable "donor" cells that already carry the two fill styles used
# throughout the sheet - green ("Supported", style 3) and red
# ("Not supported", style 5). Copying *format only* from these reuses
# the existing style indices instead of fabricating new duplicate ones.
$greenDonor = $ws.Range("B4")   # green / Supported
$redDonor   = $ws.Range("D7")   # red   / Not supported

$xlPasteFormats = -4122

function Set-SupportColor($rng, $donor) {
    $donor.Copy() | Out-Null
    $rng.PasteSpecial($xlPasteFormats) | Out-Null
}

# Row 39 (virtualbox): B39 green -> red, D39 red -> green
Set-SupportColor $ws.Range("B39") $redDonor
Set-SupportColor $ws.Range("D39") $greenDonor

# Row 40 (visual-studio-2017-community): C40 green -> red
Set-SupportColor $ws.Range("C40") $redDonor

# Row 42 (visual-studio-code): C42 red -> green
Set-SupportColor $ws.Range("C42") $greenDonor

# Row 43 (vlc): B43 red -> green, C43 green -> red
Set-SupportColor $ws.Range("B43") $greenDonor
Set-SupportColor $ws.Range("C43") $redDonor

# Row 44 (win32-disk-imager): B44 green -> red
Set-SupportColor $ws.Range("B44") $redDonor

# New row 45: virtualbox-guest-additions
$ws.Range("A45").Value2 = "virtualbox-guest-additions"
Set-SupportColor $ws.Range("B45") $greenDonor
Set-SupportColor $ws.Range("C45") $greenDonor
Set-SupportColor $ws.Range("D45") $redDonor

$ws.Application.CutCopyMode = $false

# Scroll/selection moved down to show the newly added row.
$ws.Range("A40").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 22
